$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (date, category, amount)
$data = @(
    @(45779, "Бытовые услуги", 432),
    @(45779, "Кафе, рестораны, фастфуд", 28857),
    @(45779, "Комиссия", 200),
    @(45779, "Одежда и обувь", 16469),
    @(45779, "Путешествия", 117704),
    @(45779, "Супермаркеты", 7989),
    @(45779, "Такси и каршеринг", 6584),
    @(45779, "Хобби и развлечения", 7870)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove old rows 10-15 (now unused)
$ws.Rows("10:15").Delete()
